$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily data point for 2026/01/18 (day "日") was logged, inserted as a
# new row 673; every existing row from 673 through 714 shifts down by one
# (ending at 715), exactly like Excel's native "insert row" behavior.
$ws.Rows.Item(673).Insert()

# Populate the freshly inserted row 673 with the new reading.
# Column A holds the date as plain text (matching the rest of the column,
# which stores "YYYY/MM/DD" as inline/shared strings, not real dates), so we
# temporarily force a text number format before assigning the value, then
# clear the format again so the cell ends up with the same default styling
# as its neighbours (no explicit style index).
$ws.Cells.Item(673, 1).NumberFormat = "@"
$ws.Cells.Item(673, 1).Value = "2026/01/18"
$ws.Cells.Item(673, 1).ClearFormats()

$ws.Cells.Item(673, 2).Value = "日"
$ws.Cells.Item(673, 3).Value = 7
$ws.Cells.Item(673, 4).Value = 174
